$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "RelationShip" to "Sheet1"
$ws.Name = "Sheet1"

# Row 4: same highlighted style as row 2, new component "org.emoflon.ibex.tgg.core.runtime"
$ws.Range("A2:D2").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Range("A4").Value = "org.emoflon.ibex.tgg.core.runtime"

# Set the text values first, in row order, so new shared-string entries are
# appended in the same order as the target (Extension, Extension_1, Extension_2).
$ws.Range("B5").Value = "(org.emoflon.ibex.tgg.ide).IbexTGGBuilderExtension"
$ws.Range("B6").Value = "(org.emoflon.ibex.tgg.ide).IbexTGGBuilderExtension_1"
$ws.Range("B7").Value = "(org.emoflon.ibex.tgg.ide).IbexTGGBuilderExtension_2"

# Row 6: highlighted in yellow (fill registered before the red one, matching
# the target fill table order of yellow, then red, then green)
$ws.Range("D6").Interior.Color = 65535

# Row 5: highlighted in red
$ws.Range("D5").Interior.Color = 255

# Row 7: highlighted in green
$ws.Range("D7").Interior.Color = 5287936

# Update the active selection shown when the sheet is opened
$ws.Range("K9").Select()
